$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 288
$firstDataRow = 2

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    # Column C (3) = "Förändrad" date -> bump 45184 to 45186
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value2 = 45186
    }

    # Beteckning for this row, used as the HYPERLINK friendly name
    $beteckning = $ws.Cells.Item($r, 1).Value2

    # Columns S..Y (19..25) hold HYPERLINK(...) formulas; append a friendly
    # name argument if one isn't already present.
    for ($col = 19; $col -le 25; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $f = $cell.Formula
        if ([string]::IsNullOrEmpty($f)) {
            continue
        }
        if ($f.IndexOf("HYPERLINK(") -ge 0 -and $f.TrimEnd().EndsWith(")") -and $f.IndexOf(",") -lt 0) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ", `"$beteckning`")"
            $cell.Formula = $newFormula
        }
    }
}
